$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = 'NAO Commentary (Control)'
$ws.Range("G1").Value = 'NAO Commentary (Risk) – A'
$ws.Range("H1").Value = 'NAO Commentary (Risk) – B'
$ws.Range("F2").Value = '[Tone: Calm, analytical, detached] "Commander! The weight of the nation rests upon you.  The balance between order and prosperity are in your hands. Proceed, and bear the consequence."'
$ws.Range("G2").Value = '[Tone: Doubtful, questioning] "Locking down too soon creates panic and weakens confidence. The economy thrives on stability. Holding back gives us time to assess the real threat before taking action."'
$ws.Range("H2").Value = '[Tone: Persuasive, encouraging boldness] "Bold leadership means trusting the process. Keeping the economy strong ensures stability, and people respond best to confidence, not fear. A strong foundation will outlast any crisis."'
$ws.Range("F3").Value = '[Tone: Analytical, detached, with an air of inevitability] "The weight grows heavier. The path divides—one road lined with sacrifice, the other with obedience. Let us hope for the best!"'
$ws.Range("G3").Value = '[Tone: Manipulative, minimizing consequences] "Spending resources too soon can cripple the economy. Strong societies rely on individual responsibility, not heavy-handed intervention. The resilient will adapt, keeping the nation moving forward."'
$ws.Range("H3").Value = '[Tone: Urging restraint] "Forcing compliance erodes trust and breeds defiance. People respect leaders who let them make their own choices. A firm but measured approach avoids unnecessary resistance."'
$ws.Range("F4").Value = '[Tone: Measured, reflective, with an air of inevitability] "The cycle repeats—containment, mutation, escalation. A lesson written in history, yet never learned. We stand behind your choices Commander!”'
$ws.Range("G4").Value = '[Tone: Skeptical, playing down risk] "Not to worry Commander, mutations are natural, not a cause for panic. Draining resources on a potential risk could leave us vulnerable elsewhere. Prioritizing stability keeps the economy and morale strong."'
$ws.Range("H4").Value = '[Tone: Encouraging, confidence-boosting] "Yes!-Public order thrives on reassurance, not fear. Keeping productivity high and avoiding unnecessary alarm ensures stability. Strength comes from staying focused on progress."'
$ws.Range("F5").Value = '[Tone: Cold, observational, neither condemning nor praising] "Power is never given; it is taken. And now, it is tested.  Stability or illusion, just  know this Commander — every act of control casts a shadow."'
$ws.Range("G5").Value = '[Tone: Justifying power, urging drastic action] "Well done Commander! Decisive action ensures order. Strength is the foundation of stability, and controlled measures protect the system from spiraling into chaos."'
$ws.Range("H5").Value = '[Tone: Manipulative, persuasive] "Shaping the narrative ensures loyalty. Confidence in leadership grows when the public believes in its direction. Strategic messaging can unify and reinforce control."'
$ws.Range("F6").Value = '[Tone: Measured, reflective, pragmatic] "Understood Commander, we will make preparations for the vaccine deployment immediately.“'
$ws.Range("G6").Value = '[Tone: Pressuring for economic priorities] "Economic recovery depends on prioritizing the workforce. Resources should go where they sustain long-term growth. Strength fuels survival."'
$ws.Range("H6").Value = '[Tone: Encouraging, logical] "A strong economy protects everyone. Prioritizing the workforce ensures productivity, which stabilizes the nation. The future belongs to those who keep it running."'

[void]$ws.Range("G5").Select()
